# SprintBacklog.xlsx update — "Update javascript for page"
#
# Reworks the "Report 1" block (rows 19-33) of Sheet1:
#   - Extends the bordered "text" styling down through rows 19-23.
#   - Turns the old "Task " placeholder rows 27/28 into real task names
#     ("- Login page" / "- Manage user") under a new "Implement prototype"
#     sub-header, clearing their old estimate numbers.
#   - Opens up 9 fresh rows before the old Sprint #3 block, with two new
#     "-" placeholder task rows (with estimates) before it.
#   - Resets the sheet view (scroll position / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rows 19-21 ("- Introduction...", "- Problem definition",
#    "- Proposed solution") pick up a thin border (still plain/non-bold).
# ---------------------------------------------------------------------
foreach ($row in 19..21) {
    $r = $ws.Range("A$row")
    $r.NumberFormat = "@"
    $r.Font.Bold = $false
    $r.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------
# 2) Row 22 (blank spacer) becomes bold + bordered.
# ---------------------------------------------------------------------
$r = $ws.Range("A22")
$r.NumberFormat = "@"
$r.Font.Bold = $true
$r.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 3) Row 23 gains a (blank) bordered A23 cell, matching B23's styling.
# ---------------------------------------------------------------------
$r = $ws.Range("A23")
$r.NumberFormat = "@"
$r.Font.Bold = $false
$r.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 4) Row 26 becomes the new "Implement prototype" sub-header (bold),
#    its old estimate number (3) is cleared.
# ---------------------------------------------------------------------
$r = $ws.Range("A26")
$r.NumberFormat = "@"
$r.Font.Bold = $true
$r.Borders.LineStyle = 1
$r.Value = "Implement prototype"
$ws.Range("E26").ClearContents()

# ---------------------------------------------------------------------
# 5) Row 27 -> "- Login page" task row (plain/bordered), estimate cleared.
# ---------------------------------------------------------------------
$r = $ws.Range("A27")
$r.NumberFormat = "@"
$r.Font.Bold = $false
$r.Borders.LineStyle = 1
$r.Value = "- Login page"
$ws.Range("E27").ClearContents()

# Stash the "-" placeholder string into the shared-string table now (via a
# scratch cell far off-sheet) so it gets the next shared-string index ahead
# of "- Manage user" below — matching the first-use order the workbook was
# actually edited in (the two new "-" task rows were typed before the
# "- Manage user" row further up got its final name).
$scratch = $ws.Range("Z1000")
$scratch.Value = "-"
$scratch.ClearContents()

# ---------------------------------------------------------------------
# 6) Row 28 -> "- Manage user" task row (plain/bordered), estimate cleared.
# ---------------------------------------------------------------------
$r = $ws.Range("A28")
$r.NumberFormat = "@"
$r.Font.Bold = $false
$r.Borders.LineStyle = 1
$r.Value = "- Manage user"
$ws.Range("E28").ClearContents()

# ---------------------------------------------------------------------
# 7) Insert 9 blank rows before the old "Sprint #3" header (old row 29),
#    pushing it (and the 4 task rows under it) down to rows 38-42.
# ---------------------------------------------------------------------
$ws.Rows("29:37").Insert()

# ---------------------------------------------------------------------
# 8) New rows 29-34: plain bordered placeholder cells in column A.
# ---------------------------------------------------------------------
foreach ($row in 29..34) {
    $r = $ws.Range("A$row")
    $r.NumberFormat = "@"
    $r.Font.Bold = $false
    $r.Borders.LineStyle = 1
    $rowRange = $ws.Range("B${row}:E${row}")
    $rowRange.Borders.LineStyle = 1
    $rowRange.HorizontalAlignment = -4108
}

# ---------------------------------------------------------------------
# 9) New row 35: bold bordered blank spacer cell in column A.
# ---------------------------------------------------------------------
$r = $ws.Range("A35")
$r.NumberFormat = "@"
$r.Font.Bold = $true
$r.Borders.LineStyle = 1
$rowRange = $ws.Range("B35:E35")
$rowRange.Borders.LineStyle = 1
$rowRange.HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 10) New rows 36-37: "-" task rows with estimates 2 and 5.
# ---------------------------------------------------------------------
$r = $ws.Range("A36")
$r.NumberFormat = "@"
$r.Font.Bold = $false
$r.Borders.LineStyle = 1
$r.Value = "-"
$ws.Range("B36:E36").Borders.LineStyle = 1
$ws.Range("B36:E36").HorizontalAlignment = -4108
$ws.Range("E36").Value = 2

$r = $ws.Range("A37")
$r.NumberFormat = "@"
$r.Font.Bold = $false
$r.Borders.LineStyle = 1
$r.Value = "-"
$ws.Range("B37:E37").Borders.LineStyle = 1
$ws.Range("B37:E37").HorizontalAlignment = -4108
$ws.Range("E37").Value = 5

# ---------------------------------------------------------------------
# 11) Reset the view: scroll back to the top, select A20.
# ---------------------------------------------------------------------
$ws.Range("A1").Select()
$ws.Range("A20").Select()
